# "Introduzco estilos para proximos cambios"
#
# Rewrites the data rows of Sheet1: row 2 gets new values, the old row 3
# becomes the new row 3 (with a tweaked "dia" value) and five brand-new
# rows (4-8) are appended, growing the used range from A1:K3 to A1:K8.
#
# A handful of source values look like numbers/dates ("56", "23",
# "2024-09-08", ...) but must be stored as literal text, matching the
# original file's column typing. Plain Range.Value assignment lets Excel
# auto-coerce numeric-/date-looking strings, so for those specific cells
# we temporarily force a Text number format, assign the string, then
# clear the format again so the cell ends up as plain text with no
# left-over styling (mirrors how the rest of the sheet's text cells look).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ws, $ref, [string]$value) {
    $rng = $ws.Range($ref)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

function Set-Row($ws, [int]$row, [object[]]$values) {
    $cols = @("A", "B", "C", "D", "E", "F", "G", "H", "I", "J", "K")
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ref = "$($cols[$i])$row"
        $v = $values[$i]
        if ($v -is [string]) {
            Set-TextCell $ws $ref $v
        } else {
            $ws.Range($ref).Value = $v
        }
    }
}

# id, dia, fecha, mes, unidad_numero, clase_numero, curso, caracter_clase,
# contenidos_tematicos, actividades, observaciones
Set-Row $ws 2 @(1, "56", "2024-09-08", 2, 5, 6, "6a", "GSDFGSDGDSFGD", "nnnnnn", "sssss", "ccccc")
Set-Row $ws 3 @(2, "5aaa", "2024-09-07", 23, 3, 1, "5b", "asdfa", "sdf", "sda", "fda")
Set-Row $ws 4 @(3, "BBBB", "2024-09-07", 32, 23, 23, "6a", "asdf", "asdf", "asd", "fasd")
Set-Row $ws 5 @(5, "Hola te amo mucho", "2024-09-08", "2", "5", "5", "5b", "Dhjfnsjfj", "Ndjfjdndjfjd", "Djfjdjtjdj", "Jfjfjfjdd")
Set-Row $ws 6 @(6, "aa", "2024-09-13", "fasdf", "3", "43434", "6a", "sdgsdfg", "dfg", "sdfgsd", "fgsdfgsdf")
Set-Row $ws 7 @(7, "asdfasd", "223223-02-23", "as2", "23", "232", "6a", "sadfsa", "dfasd", "fasdf", "asf")
Set-Row $ws 8 @(8, "zzzz", "2024-09-08", "23", "23", "32", "6a", "asdf", "sdaf", "sdfas", "df")
